# The commit re-orders the four species observation records that occupy
# rows 2-5 of the sheet. Comparing the old and new values shows that the
# whole rows (including the optional "Ålder-Stadium"/"Publik kommentar"
# cells that only the "Knärot" record carries) simply get rotated:
#
#   old row 2 (TaxonId 4412,   Aggvaxskivling)      -> new row 5
#   old row 3 (TaxonId 4769,   Svavelriska)          -> new row 4
#   old row 4 (TaxonId 249278, Barrviolspindling)    -> new row 2
#   old row 5 (TaxonId 220787, Knarot)               -> new row 3
#
# That is exactly the order obtained by sorting rows 2-5 in descending
# order of the "TaxonId" column (column E): 249278, 220787, 4769, 4412.
# Using Range.Sort keeps every cell of a record together (including the
# sparsely-populated helper columns), which reproduces the target layout
# precisely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:AY5")
$sortKey = $ws.Range("E2:E5")

$dataRange.Sort($sortKey, 2)
